$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.802.60"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.336.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.26%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.01"
$ws.Range("D6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +2.11%  "

$ws.Range("E9").Value = "  +2.52%  "

$ws.Range("E10").Value = "  +3.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.906.10"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.22%  "

$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.79"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.787.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.349.82"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.55%  "

$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.539"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.63"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.178"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.80"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0961"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.03%  "

$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.48"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.99"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.57"
$ws.Range("D32").ClearFormats()

$ws.Range("E33").Value = "  +6.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.67"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.90"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.91%  "

$ws.Range("E36").Value = "  +9.63%  "

$ws.Range("E37").Value = "  +11.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.27"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.847.83"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0735"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.43%  "

$ws.Range("E41").Value = "  +8.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.58%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.749"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.55%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.84%  "

$ws.Range("E45").Value = "  +2.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.377.04"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.99"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.48%  "

$ws.Range("E48").Value = "  +3.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.28"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "282.88"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.36%  "
